# Apply the edit described by the diff:
#  1. Update the fixed "date" placeholder text (slide master + every slide
#     layout) from 24/09/2020 to 07/10/2020.
#  2. Update the scenario title textbox on slide 4 from
#     "Cenário: Cancelar compra" to "Cenário: Realizar a troca do produto"
#     (its height auto-grows because the shape already has spAutoFit, which
#     matches the taller <a:ext cy="1569660"/> from the diff).

$p = $ppt.ActivePresentation
$newDate = "07/10/2020"

function Update-DatePlaceholderShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$sm = $p.Designs.Item(1).SlideMaster

# Slide master's own date placeholder.
Update-DatePlaceholderShapes $sm.Shapes

# Every slide layout hanging off the master (each has its own date placeholder).
for ($i = 1; $i -le $sm.CustomLayouts.Count; $i++) {
    $layout = $sm.CustomLayouts.Item($i)
    Update-DatePlaceholderShapes $layout.Shapes
}

# Slide 4: "Cenário: Cancelar compra" -> "Cenário: Realizar a troca do produto".
$s = $p.Slides.Item(4)
$titleBox = $null
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $shp = $s.Shapes.Item($j)
    if ($shp.Name -eq "CaixaDeTexto 2") {
        $titleBox = $shp
    }
}
if ($titleBox -eq $null) {
    $titleBox = $s.Shapes.Item(7)
}
$titleBox.TextFrame.TextRange.Text = "Cenário: Realizar a troca do produto"
